$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 1,16

$arr[0,0] = 2
$arr[0,1] = 1
$arr[0,2] = 1.736858
$arr[0,3] = 3.473716
$arr[0,4] = 0.01904277991942104
$arr[0,5] = 0.01461981882753793
$arr[0,6] = 2
$arr[0,7] = 1
$arr[0,8] = 164.5772705
$arr[0,9] = 329.154541
$arr[0,10] = 0.2320765473082729
$arr[0,11] = 0.1805598029509348
$arr[0,12] = 285.847348886089
$arr[0,13] = 1143.389395544356
$arr[0,14] = 0.004419382614850547
$arr[0,15] = 0.002639751606678615
$ws.Range("E2:T2").Value = $arr

$arr[0,0] = 2
$arr[0,1] = 1
$arr[0,2] = 1.736858
$arr[0,3] = 3.473716
$arr[0,4] = 0.01904277991942104
$arr[0,5] = 0.01461981882753793
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 118.764328
$arr[0,9] = 356.292984
$arr[0,10] = 0.1674740084210307
$arr[0,11] = 0.1954467673099505
$arr[0,12] = 206.276773201424
$arr[0,13] = 1237.660639208544
$arr[0,14] = 0.003189170684584953
$arr[0,15] = 0.00285739632849944
$ws.Range("E3:T3").Value = $arr

$arr[0,0] = 2
$arr[0,1] = 1
$arr[0,2] = 1.736858
$arr[0,3] = 3.473716
$arr[0,4] = 0.01904277991942104
$arr[0,5] = 0.01461981882753793
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 84.531957
$arr[0,9] = 253.595871
$arr[0,10] = 0.1192016653221345
$arr[0,11] = 0.1391116171686985
$arr[0,12] = 146.820005771106
$arr[0,13] = 880.9200346266359
$arr[0,14] = 0.002269931078757889
$arr[0,15] = 0.002033786639812188
$ws.Range("E4:T4").Value = $arr

$arr[0,0] = 2
$arr[0,1] = 1
$arr[0,2] = 1.736858
$arr[0,3] = 3.473716
$arr[0,4] = 0.01904277991942104
$arr[0,5] = 0.01461981882753793
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 79.63570166666666
$arr[0,9] = 238.907105
$arr[0,10] = 0.1122972730628175
$arr[0,11] = 0.1310540017808179
$arr[0,12] = 138.3159055253633
$arr[0,13] = 829.8954331521799
$arr[0,14] = 0.002138452256486362
$arr[0,15] = 0.001915985762659391
$ws.Range("E5:T5").Value = $arr

$arr[0,0] = 2
$arr[0,1] = 1
$arr[0,2] = 1.736858
$arr[0,3] = 3.473716
$arr[0,4] = 0.01904277991942104
$arr[0,5] = 0.01461981882753793
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 121.733284
$arr[0,9] = 365.199852
$arr[0,10] = 0.1716606440086599
$arr[0,11] = 0.2003326860218846
$arr[0,12] = 211.433428181672
$arr[0,13] = 1268.600569090032
$arr[0,14] = 0.003268895864682992
$arr[0,15] = 0.002928827574873994
$ws.Range("E6:T6").Value = $arr

$arr[0,0] = 2
$arr[0,1] = 1
$arr[0,2] = 1.736858
$arr[0,3] = 3.473716
$arr[0,4] = 0.01904277991942104
$arr[0,5] = 0.01461981882753793
$arr[0,6] = 2
$arr[0,7] = 1
$arr[0,8] = 139.9082645
$arr[0,9] = 279.816529
$arr[0,10] = 0.1972898618770847
$arr[0,11] = 0.1534951247677137
$arr[0,12] = 243.000788462941
$arr[0,13] = 972.0031538517638
$arr[0,14] = 0.0037569474200583
$arr[0,15] = 0.002244070915014305
$ws.Range("E7:T7").Value = $arr

$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 43.378693
$arr[0,3] = 130.136079
$arr[0,4] = 0.4756007134671516
$arr[0,5] = 0.5477033522332176
$arr[0,6] = 2
$arr[0,7] = 1
$arr[0,8] = 164.5772705
$arr[0,9] = 329.154541
$arr[0,10] = 0.2320765473082729
$arr[0,11] = 0.1805598029509348
$arr[0,12] = 7139.146891797456
$arr[0,13] = 42834.88135078474
$arr[0,14] = 0.1103757714788078
$arr[0,15] = 0.09889320935479622
$ws.Range("E8:T8").Value = $arr

$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 43.378693
$arr[0,3] = 130.136079
$arr[0,4] = 0.4756007134671516
$arr[0,5] = 0.5477033522332176
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 118.764328
$arr[0,9] = 356.292984
$arr[0,10] = 0.1674740084210307
$arr[0,11] = 0.1954467673099505
$arr[0,12] = 5151.841323663303
$arr[0,13] = 46366.57191296973
$arr[0,14] = 0.07965075789224593
$arr[0,15] = 0.1070468496388055
$ws.Range("E9:T9").Value = $arr

$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 43.378693
$arr[0,3] = 130.136079
$arr[0,4] = 0.4756007134671516
$arr[0,5] = 0.5477033522332176
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 84.531957
$arr[0,9] = 253.595871
$arr[0,10] = 0.1192016653221345
$arr[0,11] = 0.1391116171686985
$arr[0,12] = 3666.885811392201
$arr[0,13] = 33001.97230252981
$arr[0,14] = 0.05669239707367978
$arr[0,15] = 0.07619189905788021
$ws.Range("E10:T10").Value = $arr

$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 43.378693
$arr[0,3] = 130.136079
$arr[0,4] = 0.4756007134671516
$arr[0,5] = 0.5477033522332176
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 79.63570166666666
$arr[0,9] = 238.907105
$arr[0,10] = 0.1122972730628175
$arr[0,11] = 0.1310540017808179
$arr[0,12] = 3454.492654437921
$arr[0,13] = 31090.43388994129
$arr[0,14] = 0.05340866318909154
$arr[0,15] = 0.07177871609893204
$ws.Range("E11:T11").Value = $arr

$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 43.378693
$arr[0,3] = 130.136079
$arr[0,4] = 0.4756007134671516
$arr[0,5] = 0.5477033522332176
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 121.733284
$arr[0,9] = 365.199852
$arr[0,10] = 0.1716606440086599
$arr[0,11] = 0.2003326860218846
$arr[0,12] = 5280.630754517812
$arr[0,13] = 47525.6767906603
$arr[0,14] = 0.08164192476474937
$arr[0,15] = 0.1097228836960709
$ws.Range("E12:T12").Value = $arr

$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 43.378693
$arr[0,3] = 130.136079
$arr[0,4] = 0.4756007134671516
$arr[0,5] = 0.5477033522332176
$arr[0,6] = 2
$arr[0,7] = 1
$arr[0,8] = 139.9082645
$arr[0,9] = 279.816529
$arr[0,10] = 0.1972898618770847
$arr[0,11] = 0.1534951247677137
$arr[0,12] = 6069.037653908298
$arr[0,13] = 36414.22592344979
$arr[0,14] = 0.0938311990685773
$arr[0,15] = 0.08406979438673279
$ws.Range("E13:T13").Value = $arr

$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 4.494369666666667
$arr[0,3] = 13.483109
$arr[0,4] = 0.04927592954568251
$arr[0,5] = 0.05674632319163286
$arr[0,6] = 2
$arr[0,7] = 1
$arr[0,8] = 164.5772705
$arr[0,9] = 329.154541
$arr[0,10] = 0.2320765473082729
$arr[0,11] = 0.1805598029509348
$arr[0,12] = 739.6710923579948
$arr[0,13] = 4438.026554147968
$arr[0,14] = 0.01143578759436771
$arr[0,15] = 0.01024610493367129
$ws.Range("E14:T14").Value = $arr

$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 4.494369666666667
$arr[0,3] = 13.483109
$arr[0,4] = 0.04927592954568251
$arr[0,5] = 0.05674632319163286
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 118.764328
$arr[0,9] = 356.292984
$arr[0,10] = 0.1674740084210307
$arr[0,11] = 0.1954467673099505
$arr[0,12] = 533.7707932452506
$arr[0,13] = 4803.937139207255
$arr[0,14] = 0.008252437439687746
$arr[0,15] = 0.01109088542453032
$ws.Range("E15:T15").Value = $arr

$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 4.494369666666667
$arr[0,3] = 13.483109
$arr[0,4] = 0.04927592954568251
$arr[0,5] = 0.05674632319163286
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 84.531957
$arr[0,9] = 253.595871
$arr[0,10] = 0.1192016653221345
$arr[0,11] = 0.1391116171686985
$arr[0,12] = 379.917863404771
$arr[0,13] = 3419.260770642939
$arr[0,14] = 0.005873772862141524
$arr[0,15] = 0.00789407278756567
$ws.Range("E16:T16").Value = $arr

$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 4.494369666666667
$arr[0,3] = 13.483109
$arr[0,4] = 0.04927592954568251
$arr[0,5] = 0.05674632319163286
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 79.63570166666666
$arr[0,9] = 238.907105
$arr[0,10] = 0.1122972730628175
$arr[0,11] = 0.1310540017808179
$arr[0,12] = 357.9122819543828
$arr[0,13] = 3221.210537589445
$arr[0,14] = 0.005533552515615664
$arr[0,15] = 0.007436832740611121
$ws.Range("E17:T17").Value = $arr

$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 4.494369666666667
$arr[0,3] = 13.483109
$arr[0,4] = 0.04927592954568251
$arr[0,5] = 0.05674632319163286
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 121.733284
$arr[0,9] = 365.199852
$arr[0,10] = 0.1716606440086599
$arr[0,11] = 0.2003326860218846
$arr[0,12] = 547.1143790333185
$arr[0,13] = 4924.029411299867
$arr[0,14] = 0.008458737799937211
$arr[0,15] = 0.01136814334684578
$ws.Range("E18:T18").Value = $arr

$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 4.494369666666667
$arr[0,3] = 13.483109
$arr[0,4] = 0.04927592954568251
$arr[0,5] = 0.05674632319163286
$arr[0,6] = 2
$arr[0,7] = 1
$arr[0,8] = 139.9082645
$arr[0,9] = 279.816529
$arr[0,10] = 0.1972898618770847
$arr[0,11] = 0.1534951247677137
$arr[0,12] = 628.7994600847768
$arr[0,13] = 3772.796760508661
$arr[0,14] = 0.009721641333932662
$arr[0,15] = 0.008710283958408692
$ws.Range("E19:T19").Value = $arr

$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 3.071414666666667
$arr[0,3] = 9.214244
$arr[0,4] = 0.0336747584077773
$arr[0,5] = 0.03877996298854842
$arr[0,6] = 2
$arr[0,7] = 1
$arr[0,8] = 164.5772705
$arr[0,9] = 329.154541
$arr[0,10] = 0.2320765473082729
$arr[0,11] = 0.1805598029509348
$arr[0,12] = 505.4850424136674
$arr[0,13] = 3032.910254482004
$arr[0,14] = 0.007815121662717191
$arr[0,15] = 0.007002102475656847
$ws.Range("E20:T20").Value = $arr

$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 3.071414666666667
$arr[0,3] = 9.214244
$arr[0,4] = 0.0336747584077773
$arr[0,5] = 0.03877996298854842
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 118.764328
$arr[0,9] = 356.292984
$arr[0,10] = 0.1674740084210307
$arr[0,11] = 0.1954467673099505
$arr[0,12] = 364.7744988960107
$arr[0,13] = 3282.970490064096
$arr[0,14] = 0.005639646773160269
$arr[0,15] = 0.007579418402511315
$ws.Range("E21:T21").Value = $arr

$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 3.071414666666667
$arr[0,3] = 9.214244
$arr[0,4] = 0.0336747584077773
$arr[0,5] = 0.03877996298854842
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 84.531957
$arr[0,9] = 253.595871
$arr[0,10] = 0.1192016653221345
$arr[0,11] = 0.1391116171686985
$arr[0,12] = 259.632692531836
$arr[0,13] = 2336.694232786524
$arr[0,14] = 0.004014087281527604
$arr[0,15] = 0.005394743365079245
$ws.Range("E22:T22").Value = $arr

$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 3.071414666666667
$arr[0,3] = 9.214244
$arr[0,4] = 0.0336747584077773
$arr[0,5] = 0.03877996298854842
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 79.63570166666666
$arr[0,9] = 238.907105
$arr[0,10] = 0.1122972730628175
$arr[0,11] = 0.1310540017808179
$arr[0,12] = 244.5942620892911
$arr[0,13] = 2201.34835880362
$arr[0,14] = 0.003781583540242577
$arr[0,15] = 0.005082269338561276
$ws.Range("E23:T23").Value = $arr

$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 3.071414666666667
$arr[0,3] = 9.214244
$arr[0,4] = 0.0336747584077773
$arr[0,5] = 0.03877996298854842
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 121.733284
$arr[0,9] = 365.199852
$arr[0,10] = 0.1716606440086599
$arr[0,11] = 0.2003326860218846
$arr[0,12] = 373.8933938990987
$arr[0,13] = 3365.040545091888
$arr[0,14] = 0.005780630715115086
$arr[0,15] = 0.007768894149325177
$ws.Range("E24:T24").Value = $arr

$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 3.071414666666667
$arr[0,3] = 9.214244
$arr[0,4] = 0.0336747584077773
$arr[0,5] = 0.03877996298854842
$arr[0,6] = 2
$arr[0,7] = 1
$arr[0,8] = 139.9082645
$arr[0,9] = 279.816529
$arr[0,10] = 0.1972898618770847
$arr[0,11] = 0.1534951247677137
$arr[0,12] = 429.7162955731794
$arr[0,13] = 2578.297773439076
$arr[0,14] = 0.006643688435014583
$arr[0,15] = 0.005952535257414559
$ws.Range("E25:T25").Value = $arr

$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 4.242306999999999
$arr[0,3] = 12.726921
$arr[0,4] = 0.04651233350775901
$arr[0,5] = 0.05356375686797306
$arr[0,6] = 2
$arr[0,7] = 1
$arr[0,8] = 164.5772705
$arr[0,9] = 329.154541
$arr[0,10] = 0.2320765473082729
$arr[0,11] = 0.1805598029509348
$arr[0,12] = 698.1873066830434
$arr[0,13] = 4189.12384009826
$arr[0,14] = 0.0107944217677316
$arr[0,15] = 0.009671461385392996
$ws.Range("E26:T26").Value = $arr

$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 4.242306999999999
$arr[0,3] = 12.726921
$arr[0,4] = 0.04651233350775901
$arr[0,5] = 0.05356375686797306
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 118.764328
$arr[0,9] = 356.292984
$arr[0,10] = 0.1674740084210307
$arr[0,11] = 0.1954467673099505
$arr[0,12] = 503.8347400246959
$arr[0,13] = 4534.512660222264
$arr[0,14] = 0.00778960693356022
$arr[0,15] = 0.01046886312482149
$ws.Range("E27:T27").Value = $arr

$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 4.242306999999999
$arr[0,3] = 12.726921
$arr[0,4] = 0.04651233350775901
$arr[0,5] = 0.05356375686797306
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 84.531957
$arr[0,9] = 253.595871
$arr[0,10] = 0.1192016653221345
$arr[0,11] = 0.1391116171686985
$arr[0,12] = 358.610512904799
$arr[0,13] = 3227.494616143191
$arr[0,14] = 0.00554434761214339
$arr[0,15] = 0.007451340839534714
$ws.Range("E28:T28").Value = $arr

$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 4.242306999999999
$arr[0,3] = 12.726921
$arr[0,4] = 0.04651233350775901
$arr[0,5] = 0.05356375686797306
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 79.63570166666666
$arr[0,9] = 238.907105
$arr[0,10] = 0.1122972730628175
$arr[0,11] = 0.1310540017808179
$arr[0,12] = 337.8390946304116
$arr[0,13] = 3040.551851673705
$arr[0,14] = 0.005223208216709649
$arr[0,15] = 0.007019744687962638
$ws.Range("E29:T29").Value = $arr

$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 4.242306999999999
$arr[0,3] = 12.726921
$arr[0,4] = 0.04651233350775901
$arr[0,5] = 0.05356375686797306
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 121.733284
$arr[0,9] = 365.199852
$arr[0,10] = 0.1716606440086599
$arr[0,11] = 0.2003326860218846
$arr[0,12] = 516.4299628461879
$arr[0,13] = 4647.869665615691
$arr[0,14] = 0.007984337124287482
$arr[0,15] = 0.01073057128678421
$ws.Range("E30:T30").Value = $arr

$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 4.242306999999999
$arr[0,3] = 12.726921
$arr[0,4] = 0.04651233350775901
$arr[0,5] = 0.05356375686797306
$arr[0,6] = 2
$arr[0,7] = 1
$arr[0,8] = 139.9082645
$arr[0,9] = 279.816529
$arr[0,10] = 0.1972898618770847
$arr[0,11] = 0.1534951247677137
$arr[0,12] = 593.5338098462014
$arr[0,13] = 3561.202859077209
$arr[0,14] = 0.009176411853326676
$arr[0,15] = 0.008221775543477007
$ws.Range("E31:T31").Value = $arr

$arr[0,0] = 2
$arr[0,1] = 1
$arr[0,2] = 34.28457450000001
$arr[0,3] = 68.56914900000001
$arr[0,4] = 0.3758934851522086
$arr[0,5] = 0.2885867858910901
$arr[0,6] = 2
$arr[0,7] = 1
$arr[0,8] = 164.5772705
$arr[0,9] = 329.154541
$arr[0,10] = 0.2320765473082729
$arr[0,11] = 0.1805598029509348
$arr[0,12] = 5642.461691463903
$arr[0,13] = 22569.84676585561
$arr[0,14] = 0.08723606218979812
$arr[0,15] = 0.05210717319473884
$ws.Range("E32:T32").Value = $arr

$arr[0,0] = 2
$arr[0,1] = 1
$arr[0,2] = 34.28457450000001
$arr[0,3] = 68.56914900000001
$arr[0,4] = 0.3758934851522086
$arr[0,5] = 0.2885867858910901
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 118.764328
$arr[0,9] = 356.292984
$arr[0,10] = 0.1674740084210307
$arr[0,11] = 0.1954467673099505
$arr[0,12] = 4071.784451258437
$arr[0,13] = 24430.70670755062
$arr[0,14] = 0.06295238869779154
$arr[0,15] = 0.0564033543907824
$ws.Range("E33:T33").Value = $arr

$arr[0,0] = 2
$arr[0,1] = 1
$arr[0,2] = 34.28457450000001
$arr[0,3] = 68.56914900000001
$arr[0,4] = 0.3758934851522086
$arr[0,5] = 0.2885867858910901
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 84.531957
$arr[0,9] = 253.595871
$arr[0,10] = 0.1192016653221345
$arr[0,11] = 0.1391116171686985
$arr[0,12] = 2898.142177397297
$arr[0,13] = 17388.85306438378
$arr[0,14] = 0.0448071294138843
$arr[0,15] = 0.04014577447882649
$ws.Range("E34:T34").Value = $arr

$arr[0,0] = 2
$arr[0,1] = 1
$arr[0,2] = 34.28457450000001
$arr[0,3] = 68.56914900000001
$arr[0,4] = 0.3758934851522086
$arr[0,5] = 0.2885867858910901
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 79.63570166666666
$arr[0,9] = 238.907105
$arr[0,10] = 0.1122972730628175
$arr[0,11] = 0.1310540017808179
$arr[0,12] = 2730.276146650608
$arr[0,13] = 16381.65687990365
$arr[0,14] = 0.04221181334467169
$arr[0,15] = 0.03782045315209143
$ws.Range("E35:T35").Value = $arr

$arr[0,0] = 2
$arr[0,1] = 1
$arr[0,2] = 34.28457450000001
$arr[0,3] = 68.56914900000001
$arr[0,4] = 0.3758934851522086
$arr[0,5] = 0.2885867858910901
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 121.733284
$arr[0,9] = 365.199852
$arr[0,10] = 0.1716606440086599
$arr[0,11] = 0.2003326860218846
$arr[0,12] = 4173.573844427658
$arr[0,13] = 25041.44306656595
$arr[0,14] = 0.06452611773988776
$arr[0,15] = 0.05781336596798459
$ws.Range("E36:T36").Value = $arr

$arr[0,0] = 2
$arr[0,1] = 1
$arr[0,2] = 34.28457450000001
$arr[0,3] = 68.56914900000001
$arr[0,4] = 0.3758934851522086
$arr[0,5] = 0.2885867858910901
$arr[0,6] = 2
$arr[0,7] = 1
$arr[0,8] = 139.9082645
$arr[0,9] = 279.816529
$arr[0,10] = 0.1972898618770847
$arr[0,11] = 0.1534951247677137
$arr[0,12] = 4796.695317415956
$arr[0,13] = 19186.78126966382
$arr[0,14] = 0.07415997376617524
$arr[0,15] = 0.04429666470666636
$ws.Range("E37:T37").Value = $arr

